$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: update Week_Start_Date (B) and MyForecast (D) ---
$ws1.Range("B2").Value = "'2025-01-12"
$ws1.Range("B2").Style = "Normal"
$ws1.Range("D2").Value = 65
$ws1.Range("B3").Value = "'2025-01-19"
$ws1.Range("B3").Style = "Normal"
$ws1.Range("D3").Value = 61
$ws1.Range("B4").Value = "'2025-01-26"
$ws1.Range("B4").Style = "Normal"
$ws1.Range("D4").Value = 51
$ws1.Range("B5").Value = "'2025-02-02"
$ws1.Range("B5").Style = "Normal"
$ws1.Range("D5").Value = 41
$ws1.Range("B6").Value = "'2025-02-09"
$ws1.Range("B6").Style = "Normal"
$ws1.Range("D6").Value = 41
$ws1.Range("B7").Value = "'2025-02-16"
$ws1.Range("B7").Style = "Normal"
$ws1.Range("D7").Value = 49
$ws1.Range("B8").Value = "'2025-02-23"
$ws1.Range("B8").Style = "Normal"
$ws1.Range("D8").Value = 56
$ws1.Range("B9").Value = "'2025-03-02"
$ws1.Range("B9").Style = "Normal"
$ws1.Range("D9").Value = 59
$ws1.Range("B10").Value = "'2025-03-09"
$ws1.Range("B10").Style = "Normal"
$ws1.Range("D10").Value = 54
$ws1.Range("B11").Value = "'2025-03-16"
$ws1.Range("B11").Style = "Normal"
$ws1.Range("D11").Value = 49
$ws1.Range("B12").Value = "'2025-03-23"
$ws1.Range("B12").Style = "Normal"
$ws1.Range("D12").Value = 48
$ws1.Range("B13").Value = "'2025-03-30"
$ws1.Range("B13").Style = "Normal"
$ws1.Range("D13").Value = 50
$ws1.Range("B14").Value = "'2025-04-06"
$ws1.Range("B14").Style = "Normal"
$ws1.Range("D14").Value = 50
$ws1.Range("B15").Value = "'2025-04-13"
$ws1.Range("B15").Style = "Normal"
$ws1.Range("D15").Value = 49
$ws1.Range("B16").Value = "'2025-04-20"
$ws1.Range("B16").Style = "Normal"
$ws1.Range("D16").Value = 49
$ws1.Range("B17").Value = "'2025-04-27"
$ws1.Range("B17").Style = "Normal"
$ws1.Range("D17").Value = 48

# --- Summary sheet updates ---
$ws2.Range("B2").Value = "'2022-12-25 to 2025-01-05"
$ws2.Range("B2").Style = "Normal"
$ws2.Range("B4").Value = "'241"
$ws2.Range("B4").Style = "Normal"
$ws2.Range("B5").Value = "'80"
$ws2.Range("B5").Style = "Normal"
$ws2.Range("B6").Value = "'64"
$ws2.Range("B6").Style = "Normal"
$ws2.Range("B8").Value = "'8829 units"
$ws2.Range("B8").Style = "Normal"
$ws2.Range("B9").Value = "'820"
$ws2.Range("B9").Style = "Normal"
$ws2.Range("B10").Value = "'423"
$ws2.Range("B10").Style = "Normal"
$ws2.Range("B11").Value = "'218"
$ws2.Range("B11").Style = "Normal"
$ws2.Range("B12").Value = "'65"
$ws2.Range("B12").Style = "Normal"
$ws2.Range("B13").Value = "'2025-01-12"
$ws2.Range("B13").Style = "Normal"
$ws2.Range("B14").Value = "'41"
$ws2.Range("B14").Style = "Normal"
$ws2.Range("B15").Value = "'2025-02-02"
$ws2.Range("B15").Style = "Normal"

Write-Output "edits applied"
